$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B1 with the new value
$ws.Range("B1").Value = "['D20', 0.35873146904102804]"

# Remove row 2 entirely (B2 previously held "['D1', 0.7994587191534351]")
$ws.Rows("2:2").Delete()
